$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "fra" (French) rows to "spa" (Spanish) with new translations.
$ws.Range("A4").Value = "spa"
$ws.Range("C4").Value = "Extranjera"
$ws.Range("A5").Value = "spa"
$ws.Range("C5").Value = "No extranjero"

# 2. Remove the Arabic rows (rows 6 and 7) entirely.
$ws.Rows("6:7").Delete() | Out-Null

# 3. Data rows A2:C5 lose their bold styling (revert to default/normal style).
$ws.Range("A2:C5").Font.Bold = $false

# 4. Update the active selection to match the authored workbook.
$ws.Range("C5").Select() | Out-Null
